$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns G and H, reusing the same (bold/centered/
# bordered) style already used by the rest of the header row (A1:F1).
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "num_samples"
$ws.Range("H1").Value = "fractional_uncertainty"

# Fill in the new data values for num_samples (G) and fractional_uncertainty (H)
$numSamples = @(956, 731, 971, 772, 973, 783, 940, 792, 937, 804, 953, 808, 956, 796, 964, 790)
$fracUncertainty = @(0.03474455528897709, 0.03492634042692064, 0.03169069461024138, 0.032809869967679, 0.03009467527024692, 0.03067124112369637, 0.02813431984812085, 0.03630763043839624, 0.02850992436117416, 0.03292225038833987, 0.02999691922471324, 0.03267920594117515, 0.03148952966005882, 0.03016095376771777, 0.03184115097705828, 0.02937649678565651)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $numSamples[$i]
    $ws.Cells.Item($row, 8).Value = $fracUncertainty[$i]
}
